# Update the need_to_buy data table: the forecast has rolled forward by one
# day. Dates in column A are bumped by one day and columns B:F carry the
# refreshed forecast values for rows 2-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  45945, 5706.8224689694,  5863.19699644138, 7076, 8094.554817,  48.9553893529996),
    @(3,  45946, 6074.84106324167, 6168.16542054843, 3620, 9002.414979,  228.155805679448),
    @(4,  45947, 6060.3263402325,  5559.17235453722, 3620, 8978.086464,  202.37218659603),
    @(5,  45948, 2017.7268977251,  4071.78588424182, 3620, 4646.325592,  128.349357438197),
    @(6,  45949, 2050.75371187105, 4236.34188404363, 3620, 4801.504815,  140.295541132191),
    @(7,  45950, 6392.95297294923, 6647.88533799522, 3620, 9928.524642,  273.477375293583),
    @(8,  45951, 6392.95297294923, 6672.62751759494, 3620, 9928.524642,  274.508299443571),
    @(9,  45952, 6392.95297294923, 6603.50371951683, 3620, 9928.524642,  271.628141190317),
    @(10, 45953, 6392.95297294923, 6457.54257355422, 3620, 9928.524642,  265.546426775208),
    @(11, 45954, 6392.95297294923, 5800.70312523795, 3620, 9928.524642,  238.178116428697),
    @(12, 45955, 2044.14689306201, 4023.85164059411, 3620, 5157.784442,  146.562049563837),
    @(13, 45956, 1933.9611505005,  3961.23900907581, 3620, 5037.337553,  143.525642148971),
    @(14, 45957, 6757.64823424925, 6281.47342782827, 3620, 10717.730246, 275.898143315793),
    @(15, 45958, 6757.64823424925, 6545.56627637175, 3620, 10717.730246, 286.902012005104)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
